$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.908.71'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '2.796.82'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '354.36'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').Value = '109.25'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = '0.555'
$ws.Range('E7').Value = '  -2.00%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.600'
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.00'
$ws.Range('E10').Value = '  -0.42%  '
$ws.Range('E11').Value = '  +3.08%  '
$ws.Range('E12').Value = '  +3.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0840'
$ws.Range('E13').Value = '  -1.66%  '
$ws.Range('D14').Value = '7.68'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').Value = '3.234.03'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').Value = '2.807.70'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '0.931'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').Value = '51.799.97'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').Value = '7.72'
$ws.Range('E19').Value = '  +3.97%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').Value = '13.22'
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('D22').Value = '0.0₃0968'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('D23').Value = '70.17'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '267.81'
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('D26').Value = '26.23'
$ws.Range('E26').Value = '  -1.77%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '0.162'
$ws.Range('E28').Value = '  +11.67%  '
$ws.Range('D29').Value = '10.28'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('D30').Value = '36.76'
$ws.Range('E30').Value = '  +8.20%  '
$ws.Range('D31').Value = '6.24'
$ws.Range('E31').Value = '  +9.07%  '
$ws.Range('B32').Value = 'OKB'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D32').Value = '52.06'
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('B33').Value = 'VeChain'
$ws.Range('C33').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D33').Value = '0.0454'
$ws.Range('E33').Value = '  -2.82%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').Value = '5.57'
$ws.Range('E34').Value = '  +6.12%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').Value = '1.92'
$ws.Range('E35').Value = '  -14.70%  '
$ws.Range('E36').Value = '  -1.44%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = '18.63'
$ws.Range('E38').Value = '  +2.92%  '
$ws.Range('D39').Value = '3.17'
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('D40').Value = '1.98'
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('D41').Value = '2.57'
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').Value = '0.115'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('D43').Value = '121.01'
$ws.Range('E43').Value = '  -1.09%  '
$ws.Range('D44').Value = '22.25'
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('E45').Value = '  -2.32%  '
$ws.Range('D46').Value = '2.135.16'
$ws.Range('E46').Value = '  +2.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.30'
$ws.Range('E47').Value = '  +1.57%  '
$ws.Range('D48').Value = '2.31'
$ws.Range('E48').Value = '  +6.08%  '
$ws.Range('D49').Value = '5.47'
$ws.Range('E49').Value = '  -4.15%  '
$ws.Range('D50').Value = '0.914'
$ws.Range('E50').Value = '  -2.52%  '
$ws.Range('D51').Value = '1.33'
$ws.Range('E51').Value = '  +8.91%  '
